# Insert a new weekly price record as row 36, shifting the existing
# rows 36:67 down to 37:68 (dimension grows from A1:R67 to A1:R68).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("36:36").Insert()

$ws.Cells.Item(36, 1).Value = 1
$ws.Cells.Item(36, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(36, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(36, 4).Value = 44673
$ws.Cells.Item(36, 5).Value = 15
$ws.Cells.Item(36, 6).Value = 100112021
$ws.Cells.Item(36, 7).Value = "Ají"
$ws.Cells.Item(36, 8).Value = "Inferno"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 120
$ws.Cells.Item(36, 11).Value = 25000
$ws.Cells.Item(36, 12).Value = 26000
$ws.Cells.Item(36, 13).Value = 25500
$ws.Cells.Item(36, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(36, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(36, 16).Value = 1700
$ws.Cells.Item(36, 17).Value = 15
$ws.Cells.Item(36, 18).Value = "Hortaliza"
